# Apply updated real-effort ranking values (female workers sheet)
# See commit message: "include no rank decision in binary"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F (realeffort) updates for rows 2-13
$ws.Range("F2").Value = 7.360079283446961
$ws.Range("F3").Value = 6.358750456454161
$ws.Range("F4").Value = 6.017579664918089
$ws.Range("F5").Value = 5.309718579672998
$ws.Range("F6").Value = 5.244195657518464
$ws.Range("F7").Value = 4.045026469112039
$ws.Range("F8").Value = 1.242073243576292
$ws.Range("F9").Value = 1.014010395470444
$ws.Range("F10").Value = 0.4834459824271087
$ws.Range("F11").Value = 0.1753989618967279
$ws.Range("F12").Value = 0.1348575153764161
$ws.Range("F13").Value = 0.1140016948445168

# Rows 10-13: prolificid (B numeric id, C prolific hash), name (D), and race (G) were
# reshuffled among these four rows while re_rank (H) stays sequential 9..12
$ws.Range("B10").Value = 30
$ws.Range("C10").Value = "60d5775a99b502eec8cf56b4"
$ws.Range("D10").Value = "Shadaisia"
$ws.Range("G10").Value = "Black or African American"

$ws.Range("B11").Value = 32
$ws.Range("C11").Value = "6036f9b3b1842f8b659b18c7"
$ws.Range("D11").Value = "Kellie"
$ws.Range("G11").Value = "White"

$ws.Range("B12").Value = 33
$ws.Range("C12").Value = "60cb36ee9f58331a33cf5506"
$ws.Range("D12").Value = "Shaniek"
$ws.Range("G12").Value = "Black or African American"

$ws.Range("B13").Value = 21
$ws.Range("C13").Value = "5c0e89c6c323400001e6c4a5"
$ws.Range("D13").Value = "Bri"
$ws.Range("G13").Value = "Black or African American"
